$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.359677419354839
    "C2" = 0.298546895640687
    "D2" = 0.451197053406998
    "E2" = 0.402173913043478
    "F2" = 0.304578130911843

    "B3" = 0.435483870967742
    "C3" = 0.416116248348745
    "D3" = 0.532228360957643
    "E3" = 0.414596273291925
    "F3" = 0.335981838819523

    "B4" = 0.401612903225806
    "C4" = 0.498018494055482
    "D4" = 0.530386740331492
    "E4" = 0.571428571428571
    "F4" = 0.386681800983731

    "B5" = 0.543548387096774
    "C5" = 0.535006605019815
    "D5" = 0.616942909760589
    "E5" = 0.503105590062112
    "F5" = 0.416193719258418

    "B6" = 0.72741935483871
    "C6" = 0.684280052840158
    "D6" = 0.74585635359116
    "E6" = 0.698757763975155
    "F6" = 0.437381763147938

    "B7" = 0.401612903225806
    "C7" = 0.467635402906209
    "D7" = 0.441988950276243
    "E7" = 0.498447204968944
    "F7" = 0.342792281498297
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
